$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 1 (kept its original style s="1") with the surviving record's values
$ws.Range("A1").Value = "nen"
$ws.Range("B1").Value = "ti"
$ws.Range("C1").Value = "24/10/2025"
$ws.Range("D1").Value = "17:29"
$ws.Range("E1").Value = "Present"

# Delete the remaining rows (old header data row duplicates / deleted record rows)
$ws.Range("A2:E6").EntireRow.Delete()
